$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 119.0815153333333
$ws.Range("H2").Value = 357.244546
$ws.Range("I2").Value = 0.431812569872284
$ws.Range("J2").Value = 0.4318125698722839
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 328.6107208428464
$ws.Range("R2").Value = 2957.496487585618
$ws.Range("S2").Value = 0.1111514631934493
$ws.Range("T2").Value = 0.1111514631934493

$ws.Range("G3").Value = 119.0815153333333
$ws.Range("H3").Value = 357.244546
$ws.Range("I3").Value = 0.431812569872284
$ws.Range("J3").Value = 0.4318125698722839
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 862.7726100939807
$ws.Range("R3").Value = 7764.953490845826
$ws.Range("S3").Value = 0.2918299128196716
$ws.Range("T3").Value = 0.2918299128196716

$ws.Range("G4").Value = 119.0815153333333
$ws.Range("H4").Value = 357.244546
$ws.Range("I4").Value = 0.431812569872284
$ws.Range("J4").Value = 0.4318125698722839
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 61.92611919414713
$ws.Range("R4").Value = 557.335072747324
$ws.Range("S4").Value = 0.02094630004969676
$ws.Range("T4").Value = 0.02094630004969675

$ws.Range("G5").Value = 119.0815153333333
$ws.Range("H5").Value = 357.244546
$ws.Range("I5").Value = 0.431812569872284
$ws.Range("J5").Value = 0.4318125698722839
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 23.31107989094578
$ws.Range("R5").Value = 209.799719018512
$ws.Range("S5").Value = 0.007884893809466297
$ws.Range("T5").Value = 0.007884893809466295

$ws.Range("I6").Value = 0.4460879372303943
$ws.Range("J6").Value = 0.4460879372303942
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 339.4743202031719
$ws.Range("R6").Value = 3055.268881828548
$ws.Range("S6").Value = 0.1148260388778655
$ws.Range("T6").Value = 0.1148260388778655

$ws.Range("I7").Value = 0.4460879372303943
$ws.Range("J7").Value = 0.4460879372303942
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.3014775690072121
$ws.Range("T7").Value = 0.3014775690072121

$ws.Range("I8").Value = 0.4460879372303943
$ws.Range("J8").Value = 0.4460879372303942
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 63.97334561189601
$ws.Range("R8").Value = 575.760110507064
$ws.Range("S8").Value = 0.02163876745075241
$ws.Range("T8").Value = 0.02163876745075241

$ws.Range("I9").Value = 0.4460879372303943
$ws.Range("J9").Value = 0.4460879372303942
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 24.081724962848
$ws.Range("R9").Value = 216.735524665632
$ws.Range("S9").Value = 0.008145561894564219
$ws.Range("T9").Value = 0.008145561894564219

$ws.Range("G10").Value = 33.50679633333333
$ws.Range("H10").Value = 100.520389
$ws.Range("I10").Value = 0.1215021138451521
$ws.Range("J10").Value = 0.121502113845152
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 92.4634899498041
$ws.Range("R10").Value = 832.1714095482369
$ws.Range("S10").Value = 0.03127546226590877
$ws.Range("T10").Value = 0.03127546226590876

$ws.Range("G11").Value = 33.50679633333333
$ws.Range("H11").Value = 100.520389
$ws.Range("I11").Value = 0.1215021138451521
$ws.Range("J11").Value = 0.121502113845152
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 242.7643454777677
$ws.Range("R11").Value = 2184.879109299909
$ws.Range("S11").Value = 0.08211421752109681
$ws.Range("T11").Value = 0.08211421752109679

$ws.Range("G12").Value = 33.50679633333333
$ws.Range("H12").Value = 100.520389
$ws.Range("I12").Value = 0.1215021138451521
$ws.Range("J12").Value = 0.121502113845152
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 17.42458397295178
$ws.Range("R12").Value = 156.821255756566
$ws.Range("S12").Value = 0.005893806505043851
$ws.Range("T12").Value = 0.005893806505043848

$ws.Range("G13").Value = 33.50679633333333
$ws.Range("H13").Value = 100.520389
$ws.Range("I13").Value = 0.1215021138451521
$ws.Range("J13").Value = 0.121502113845152
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 6.559201098756445
$ws.Range("R13").Value = 59.032809888808
$ws.Range("S13").Value = 0.002218627553102641
$ws.Range("T13").Value = 0.002218627553102641

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.16474
$ws.Range("H14").Value = 0.49422
$ws.Range("I14").Value = 0.000597379052169715
$ws.Range("J14").Value = 0.000597379052169715
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 0.4546073334733333
$ws.Range("R14").Value = 4.09146600126
$ws.Range("S14").Value = 0.0001537693906164393
$ws.Range("T14").Value = 0.0001537693906164393

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.16474
$ws.Range("H15").Value = 0.49422
$ws.Range("I15").Value = 0.000597379052169715
$ws.Range("J15").Value = 0.000597379052169715
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 1.19357869598
$ws.Range("R15").Value = 10.74220826382
$ws.Range("S15").Value = 0.0004037239508024234
$ws.Range("T15").Value = 0.0004037239508024234

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.16474
$ws.Range("H16").Value = 0.49422
$ws.Range("I16").Value = 0.000597379052169715
$ws.Range("J16").Value = 0.000597379052169715
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 0.08566996185333334
$ws.Range("R16").Value = 0.77102965668
$ws.Range("S16").Value = 0.00002897757439958546
$ws.Range("T16").Value = 0.00002897757439958545

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.16474
$ws.Range("H17").Value = 0.49422
$ws.Range("I17").Value = 0.000597379052169715
$ws.Range("J17").Value = 0.000597379052169715
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 0.03224906309333334
$ws.Range("R17").Value = 0.29024156784
$ws.Range("S17").Value = 0.00001090813635126688
$ws.Range("T17").Value = 0.00001090813635126688
